$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 loses its custom style (s="2" -> default / no style)
$ws.Range("A2:K2").Style = "Normal"

# New row 3 data
$row3 = @(
    " Iberdrola.",
    " Social.",
    " People Centric.",
    " Impulsar la oferta cultural en el mundo rural.",
    " No.",
    " No.",
    " Personas del mundo rural.",
    " Desarrollo Comunitario, Medioambiental.",
    " ODS 11: Ciudades y Comunidades Sostenibles.",
    " Vecinos, nios, mayores, veraneantes y todos los habitantes del mundo rural.",
    " Noticia"
)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# New row 4 data
$row4 = @(
    " Engie.",
    " Medioambiental.",
    " Planet Positive.",
    " Descarbonizar la industria y mejorar la dependencia energética del exterior.",
    " ENAGS Renovable, Fivet Hydrogen Pontegadea, Navantia.",
    " Sí, hay otras empresas que colaboran con la empresa principal. Las empresas son: Ardian (privada), Fivet Hydrogen Pontegadea (privada), Navantia (pública) y ENAGS (privada).",
    " La industria de la región de Murcia.",
    " Medioambiente.",
    " ODS 7.",
    " La región de Murcia, Repsol, ENAGS Renovable, Engie, la industria del Valle de Escombreras y la Asociación Sectorial del Hidrógeno Verde en la Región de Murcia.",
    " Noticia"
)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}
